# "Added update end date feature"
#
# Template sheet "Rz" holds a small key/value header block:
#   B2 Auftragsnummer, B3 Date_Now, B4 Time_Now, B5 Start_Date, B6 End_Date
#
# This edit bumps the order number, rolls Start_Date back, and - the new
# "update end date" feature - turns End_Date from a free-typed text string
# into a real date value (formatted the same way as the other date cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rz")

# Auftragsnummer: 2210025 -> 2210024. The cell had a one-off "general,
# right aligned" style; drop it back to the sheet default now that it's
# just a plain number again.
$ws.Range("B2").Value = 2210024
$ws.Range("B2").ClearFormats()

# Start_Date: 25-Nov-2022 -> 25-Oct-2022 (serial 44890 -> 44859). Keep the
# existing date formatting on the cell as-is.
$ws.Range("B5").Value = 44859

# End_Date: used to be the literal text "15-02-2023"; the new feature
# stores a proper date value instead (20-Jan-2023, serial 44946), using
# the same date format already applied to the cell.
$ws.Range("B6").Value = 44946

# Reflect where the user's cursor ended up after editing End_Date.
$ws.Range("B6").Select()
